# "env values created + update of all tables according to new style"
#
# Row 10 ("Leasehold buildings") gets its depreciation-life columns filled in
# with the new environment values:
#   - B10 ("Over the period of lease") loses the old ad-hoc alignment style,
#     falling back to the sheet's default/Normal formatting (matching the
#     refreshed table style used elsewhere on the sheet).
#   - C10, previously blank, is given the placeholder value "." with the
#     same default formatting.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Drop the leftover alignment-only style from B10 so it matches the rest of
# the table's default formatting.
$ws.Range("B10").Style = "Normal"

# C10 becomes a real (default-styled) cell holding the new placeholder value.
$ws.Range("C10").Style = "Normal"
$ws.Range("C10").Value = "."

# Leave the workbook's selection on C10, the last cell touched.
$ws.Range("C10").Select() | Out-Null
